# Auto-generated edit script applying odds updates per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 6).Value2 = 2.94  # F2: 2.86 -> 2.94
$ws.Cells.Item(2, 7).Value2 = 2.98  # G2: 2.88 -> 2.98
$ws.Cells.Item(2, 8).Value2 = 2.56  # H2: 2.62 -> 2.56
$ws.Cells.Item(2, 9).Value2 = 2.6  # I2: 2.66 -> 2.6
$ws.Cells.Item(2, 10).Value2 = 3.6  # J2: 3.65 -> 3.6
$ws.Cells.Item(2, 14).Value2 = 4.7  # N2: 4.8 -> 4.7
$ws.Cells.Item(2, 15).Value2 = 1.25  # O2: 1.24 -> 1.25
$ws.Cells.Item(2, 17).Value2 = 1.76  # Q2: 1.75 -> 1.76
$ws.Cells.Item(2, 18).Value2 = 1.48  # R2: 1.5 -> 1.48
$ws.Cells.Item(2, 19).Value2 = 2.92  # S2: 2.88 -> 2.92
$ws.Cells.Item(2, 20).Value2 = 1.63  # T2: 1.62 -> 1.63
$ws.Cells.Item(2, 21).Value2 = 2.5  # U2: 2.52 -> 2.5
$ws.Cells.Item(2, 22).Value2 = 1.62  # V2: 1.6 -> 1.62
$ws.Cells.Item(2, 23).Value2 = 1.5  # W2: 1.53 -> 1.5
$ws.Cells.Item(2, 25).Value2 = 13.5  # Y2: 14 -> 13.5
$ws.Cells.Item(2, 27).Value2 = 36  # AA2: 38 -> 36
$ws.Cells.Item(2, 31).Value2 = 24  # AE2: 25 -> 24
$ws.Cells.Item(2, 33).Value2 = 13  # AG2: 12.5 -> 13
$ws.Cells.Item(2, 34).Value2 = 15.5  # AH2: 15 -> 15.5
$ws.Cells.Item(2, 36).Value2 = 46  # AJ2: 44 -> 46
$ws.Cells.Item(2, 37).Value2 = 29  # AK2: 28 -> 29
$ws.Cells.Item(2, 41).Value2 = 17  # AO2: 18 -> 17
# Row 4
$ws.Cells.Item(4, 6).Value2 = 1.5  # F4: 1.04 -> 1.5
$ws.Cells.Item(4, 7).Value2 = 990  # G4: 1000 -> 990
$ws.Cells.Item(4, 9).Value2 = 12  # I4: 990 -> 12
$ws.Cells.Item(4, 10).Value2 = 1.01  # J4: 1.03 -> 1.01
$ws.Cells.Item(4, 17).Value2 = 1.31  # Q4: 1.3 -> 1.31
$ws.Cells.Item(4, 19).Value2 = 1.31  # S4: 1.3 -> 1.31
$ws.Cells.Item(4, 22).Value2 = 1.09  # V4: 1.01 -> 1.09
# Row 5
$ws.Cells.Item(5, 12).Value2 = 1.32  # L5: 1.37 -> 1.32
$ws.Cells.Item(5, 22).Value2 = 1.56  # V5: 1.54 -> 1.56
$ws.Cells.Item(5, 28).Value2 = 15  # AB5: 14.5 -> 15
$ws.Cells.Item(5, 29).Value2 = 10  # AC5: 9.800000000000001 -> 10
$ws.Cells.Item(5, 30).Value2 = 15  # AD5: 14.5 -> 15
$ws.Cells.Item(5, 31).Value2 = 36  # AE5: 980 -> 36
$ws.Cells.Item(5, 32).Value2 = 23  # AF5: 21 -> 23
$ws.Cells.Item(5, 41).Value2 = 28  # AO5: 29 -> 28
# Row 6
$ws.Cells.Item(6, 7).Value2 = 4.7  # G6: 4.8 -> 4.7
$ws.Cells.Item(6, 10).Value2 = 3.8  # J6: 3.75 -> 3.8
$ws.Cells.Item(6, 14).Value2 = 4.3  # N6: 4.4 -> 4.3
$ws.Cells.Item(6, 15).Value2 = 1.21  # O6: 1.23 -> 1.21
$ws.Cells.Item(6, 19).Value2 = 2.78  # S6: 2.76 -> 2.78
$ws.Cells.Item(6, 20).Value2 = 1.73  # T6: 1.67 -> 1.73
$ws.Cells.Item(6, 22).Value2 = 2.02  # V6: 2 -> 2.02
$ws.Cells.Item(6, 23).Value2 = 1.27  # W6: 1.26 -> 1.27
$ws.Cells.Item(6, 25).Value2 = 13  # Y6: 13.5 -> 13
# Row 7
$ws.Cells.Item(7, 7).Value2 = 3.75  # G7: 3.8 -> 3.75
$ws.Cells.Item(7, 8).Value2 = 2.18  # H7: 2.16 -> 2.18
$ws.Cells.Item(7, 9).Value2 = 2.2  # I7: 2.18 -> 2.2
$ws.Cells.Item(7, 12).Value2 = 1.38  # L7: 1.39 -> 1.38
$ws.Cells.Item(7, 16).Value2 = 2.16  # P7: 2.18 -> 2.16
$ws.Cells.Item(7, 17).Value2 = 1.85  # Q7: 1.83 -> 1.85
$ws.Cells.Item(7, 22).Value2 = 1.83  # V7: 1.84 -> 1.83
$ws.Cells.Item(7, 24).Value2 = 16  # X7: 16.5 -> 16
$ws.Cells.Item(7, 32).Value2 = 26  # AF7: 27 -> 26
$ws.Cells.Item(7, 34).Value2 = 16.5  # AH7: 16 -> 16.5
$ws.Cells.Item(7, 36).Value2 = 65  # AJ7: 70 -> 65
$ws.Cells.Item(7, 37).Value2 = 38  # AK7: 40 -> 38
$ws.Cells.Item(7, 41).Value2 = 14  # AO7: 13.5 -> 14
# Row 8
$ws.Cells.Item(8, 8).Value2 = 2.06  # H8: 2.2 -> 2.06
$ws.Cells.Item(8, 19).Value2 = 3.75  # S8: 4.2 -> 3.75
$ws.Cells.Item(8, 20).Value2 = 1.9  # T8: 1.89 -> 1.9
$ws.Cells.Item(8, 31).Value2 = 36  # AE8: 980 -> 36
# Row 9
$ws.Cells.Item(9, 20).Value2 = 1.59  # T9: 1.58 -> 1.59
$ws.Cells.Item(9, 28).Value2 = 970  # AB9: 1000 -> 970
# Row 10
$ws.Cells.Item(10, 40).Value2 = 23  # AN10: 22 -> 23
# Row 11
$ws.Cells.Item(11, 6).Value2 = 3.2  # F11: 3.15 -> 3.2
$ws.Cells.Item(11, 9).Value2 = 2.7  # I11: 2.72 -> 2.7
$ws.Cells.Item(11, 12).Value2 = 1.52  # L11: 1.43 -> 1.52
$ws.Cells.Item(11, 21).Value2 = 1.9  # U11: 1.89 -> 1.9
$ws.Cells.Item(11, 28).Value2 = 13  # AB11: 12.5 -> 13
$ws.Cells.Item(11, 31).Value2 = 42  # AE11: 980 -> 42
$ws.Cells.Item(11, 37).Value2 = 60  # AK11: 1000 -> 60
$ws.Cells.Item(11, 41).Value2 = 42  # AO11: 1000 -> 42
# Row 12
$ws.Cells.Item(12, 7).Value2 = 1.64  # G12: 1.68 -> 1.64
$ws.Cells.Item(12, 11).Value2 = 4.8  # K12: 4.9 -> 4.8
$ws.Cells.Item(12, 19).Value2 = 2.56  # S12: 2.64 -> 2.56
$ws.Cells.Item(12, 23).Value2 = 2.56  # W12: 2.46 -> 2.56
$ws.Cells.Item(12, 24).Value2 = 26  # X12: 1000 -> 26
# Row 13
$ws.Cells.Item(13, 7).Value2 = 1.87  # G13: 1.88 -> 1.87
$ws.Cells.Item(13, 12).Value2 = 1.3  # L13: 1.37 -> 1.3
$ws.Cells.Item(13, 20).Value2 = 1.76  # T13: 1.78 -> 1.76
$ws.Cells.Item(13, 23).Value2 = 2.14  # W13: 2.12 -> 2.14
$ws.Cells.Item(13, 31).Value2 = 75  # AE13: 80 -> 75
# Row 14
$ws.Cells.Item(14, 21).Value2 = 2.14  # U14: 2.12 -> 2.14
$ws.Cells.Item(14, 28).Value2 = 13  # AB14: 12.5 -> 13
$ws.Cells.Item(14, 31).Value2 = 970  # AE14: 46 -> 970
$ws.Cells.Item(14, 37).Value2 = 970  # AK14: 32 -> 970
# Row 15
$ws.Cells.Item(15, 6).Value2 = 2.48  # F15: 2.5 -> 2.48
$ws.Cells.Item(15, 8).Value2 = 2.82  # H15: 2.86 -> 2.82
$ws.Cells.Item(15, 9).Value2 = 3.1  # I15: 3.2 -> 3.1
$ws.Cells.Item(15, 10).Value2 = 3.5  # J15: 3.4 -> 3.5
$ws.Cells.Item(15, 11).Value2 = 3.65  # K15: 3.75 -> 3.65
$ws.Cells.Item(15, 14).Value2 = 3.6  # N15: 3.55 -> 3.6
$ws.Cells.Item(15, 22).Value2 = 1.48  # V15: 1.46 -> 1.48
$ws.Cells.Item(15, 23).Value2 = 1.58  # W15: 1.57 -> 1.58
$ws.Cells.Item(15, 25).Value2 = 14  # Y15: 13 -> 14
$ws.Cells.Item(15, 34).Value2 = 19.5  # AH15: 18 -> 19.5
# Row 16
$ws.Cells.Item(16, 6).Value2 = 3.05  # F16: 3 -> 3.05
$ws.Cells.Item(16, 10).Value2 = 3.55  # J16: 3.5 -> 3.55
$ws.Cells.Item(16, 11).Value2 = 3.75  # K16: 3.8 -> 3.75
$ws.Cells.Item(16, 19).Value2 = 3.15  # S16: 2.84 -> 3.15
$ws.Cells.Item(16, 22).Value2 = 1.68  # V16: 1.67 -> 1.68
# Row 17
$ws.Cells.Item(17, 6).Value2 = 2  # F17: 1.9 -> 2
$ws.Cells.Item(17, 7).Value2 = 2.24  # G17: 2.12 -> 2.24
$ws.Cells.Item(17, 8).Value2 = 3.55  # H17: 3.6 -> 3.55
$ws.Cells.Item(17, 9).Value2 = 4.4  # I17: 5.4 -> 4.4
$ws.Cells.Item(17, 11).Value2 = 4.5  # K17: 4.7 -> 4.5
$ws.Cells.Item(17, 12).Value2 = 1.32  # L17: 1.31 -> 1.32
$ws.Cells.Item(17, 14).Value2 = 3.3  # N17: 3.35 -> 3.3
$ws.Cells.Item(17, 15).Value2 = 1.29  # O17: 1.28 -> 1.29
$ws.Cells.Item(17, 17).Value2 = 1.73  # Q17: 1.72 -> 1.73
$ws.Cells.Item(17, 19).Value2 = 2.86  # S17: 2.82 -> 2.86
$ws.Cells.Item(17, 21).Value2 = 2.1  # U17: 2.08 -> 2.1
$ws.Cells.Item(17, 22).Value2 = 1.3  # V17: 1.27 -> 1.3
$ws.Cells.Item(17, 23).Value2 = 1.8  # W17: 1.89 -> 1.8
# Row 18
$ws.Cells.Item(18, 6).Value2 = 2.2  # F18: 2.16 -> 2.2
$ws.Cells.Item(18, 7).Value2 = 2.22  # G18: 2.2 -> 2.22
$ws.Cells.Item(18, 8).Value2 = 3.6  # H18: 3.65 -> 3.6
$ws.Cells.Item(18, 9).Value2 = 3.65  # I18: 3.7 -> 3.65
$ws.Cells.Item(18, 18).Value2 = 1.48  # R18: 1.47 -> 1.48
$ws.Cells.Item(18, 21).Value2 = 2.38  # U18: 2.4 -> 2.38
$ws.Cells.Item(18, 23).Value2 = 1.82  # W18: 1.83 -> 1.82
$ws.Cells.Item(18, 25).Value2 = 16  # Y18: 16.5 -> 16
$ws.Cells.Item(18, 27).Value2 = 70  # AA18: 65 -> 70
$ws.Cells.Item(18, 29).Value2 = 8.199999999999999  # AC18: 8.4 -> 8.199999999999999
$ws.Cells.Item(18, 30).Value2 = 15  # AD18: 14.5 -> 15
$ws.Cells.Item(18, 32).Value2 = 14.5  # AF18: 14 -> 14.5
$ws.Cells.Item(18, 36).Value2 = 27  # AJ18: 26 -> 27
$ws.Cells.Item(18, 37).Value2 = 21  # AK18: 20 -> 21
$ws.Cells.Item(18, 41).Value2 = 30  # AO18: 32 -> 30
# Row 19
$ws.Cells.Item(19, 8).Value2 = 1.57  # H19: 1.59 -> 1.57
$ws.Cells.Item(19, 9).Value2 = 1.58  # I19: 1.6 -> 1.58
$ws.Cells.Item(19, 10).Value2 = 4.8  # J19: 4.7 -> 4.8
$ws.Cells.Item(19, 11).Value2 = 5  # K19: 4.8 -> 5
$ws.Cells.Item(19, 17).Value2 = 1.54  # Q19: 1.55 -> 1.54
$ws.Cells.Item(19, 18).Value2 = 1.73  # R19: 1.72 -> 1.73
$ws.Cells.Item(19, 19).Value2 = 2.34  # S19: 2.36 -> 2.34
$ws.Cells.Item(19, 22).Value2 = 2.72  # V19: 2.66 -> 2.72
$ws.Cells.Item(19, 33).Value2 = 23  # AG19: 22 -> 23
$ws.Cells.Item(19, 34).Value2 = 18  # AH19: 17.5 -> 18
$ws.Cells.Item(19, 36).Value2 = 150  # AJ19: 140 -> 150
# Row 20
$ws.Cells.Item(20, 6).Value2 = 4  # F20: 3.95 -> 4
$ws.Cells.Item(20, 8).Value2 = 1.94  # H20: 1.96 -> 1.94
$ws.Cells.Item(20, 9).Value2 = 1.95  # I20: 1.97 -> 1.95
$ws.Cells.Item(20, 16).Value2 = 2.44  # P20: 2.42 -> 2.44
$ws.Cells.Item(20, 17).Value2 = 1.68  # Q20: 1.69 -> 1.68
$ws.Cells.Item(20, 18).Value2 = 1.57  # R20: 1.56 -> 1.57
$ws.Cells.Item(20, 19).Value2 = 2.7  # S20: 2.72 -> 2.7
$ws.Cells.Item(20, 20).Value2 = 1.64  # T20: 1.63 -> 1.64
$ws.Cells.Item(20, 22).Value2 = 2.04  # V20: 2.02 -> 2.04
$ws.Cells.Item(20, 26).Value2 = 13.5  # Z20: 14 -> 13.5
$ws.Cells.Item(20, 29).Value2 = 9.199999999999999  # AC20: 9 -> 9.199999999999999
$ws.Cells.Item(20, 40).Value2 = 34  # AN20: 32 -> 34
$ws.Cells.Item(20, 41).Value2 = 9.4  # AO20: 9.6 -> 9.4
# Row 21
$ws.Cells.Item(21, 6).Value2 = 2.54  # F21: 2.46 -> 2.54
$ws.Cells.Item(21, 7).Value2 = 2.58  # G21: 2.62 -> 2.58
$ws.Cells.Item(21, 9).Value2 = 3.1  # I21: 3.2 -> 3.1
$ws.Cells.Item(21, 16).Value2 = 2.38  # P21: 2.36 -> 2.38
$ws.Cells.Item(21, 17).Value2 = 1.71  # Q21: 1.72 -> 1.71
$ws.Cells.Item(21, 19).Value2 = 1.71  # S21: 1.73 -> 1.71
$ws.Cells.Item(21, 22).Value2 = 1.48  # V21: 1.46 -> 1.48
$ws.Cells.Item(21, 23).Value2 = 1.63  # W21: 1.62 -> 1.63
$ws.Cells.Item(21, 38).Value2 = 34  # AL21: 36 -> 34
# Row 22
$ws.Cells.Item(22, 6).Value2 = 1.93  # F22: 1.91 -> 1.93
$ws.Cells.Item(22, 7).Value2 = 2.04  # G22: 2.02 -> 2.04
$ws.Cells.Item(22, 8).Value2 = 4.2  # H22: 4.8 -> 4.2
$ws.Cells.Item(22, 11).Value2 = 3.45  # K22: 3.55 -> 3.45
$ws.Cells.Item(22, 12).Value2 = 1.01  # L22: 1.54 -> 1.01
$ws.Cells.Item(22, 17).Value2 = 2.44  # Q22: 2.28 -> 2.44
$ws.Cells.Item(22, 23).Value2 = 1.96  # W22: 1.98 -> 1.96
$ws.Cells.Item(22, 30).Value2 = 22  # AD22: 24 -> 22
$ws.Cells.Item(22, 39).Value2 = 220  # AM22: 230 -> 220
# Row 23
$ws.Cells.Item(23, 6).Value2 = 3.6  # F23: 3.5 -> 3.6
$ws.Cells.Item(23, 8).Value2 = 2.22  # H23: 2.16 -> 2.22
$ws.Cells.Item(23, 9).Value2 = 2.38  # I23: 2.44 -> 2.38
$ws.Cells.Item(23, 11).Value2 = 3.45  # K23: 3.5 -> 3.45
$ws.Cells.Item(23, 12).Value2 = 1.45  # L23: 1.52 -> 1.45
$ws.Cells.Item(23, 17).Value2 = 2.4  # Q23: 2.36 -> 2.4
$ws.Cells.Item(23, 20).Value2 = 2.04  # T23: 2.06 -> 2.04
$ws.Cells.Item(23, 21).Value2 = 1.78  # U23: 1.8 -> 1.78
$ws.Cells.Item(23, 22).Value2 = 1.73  # V23: 1.69 -> 1.73
$ws.Cells.Item(23, 24).Value2 = 1000  # X23: 12 -> 1000
# Row 24
$ws.Cells.Item(24, 7).Value2 = 2.42  # G24: 2.48 -> 2.42
$ws.Cells.Item(24, 9).Value2 = 4.1  # I24: 4.2 -> 4.1
$ws.Cells.Item(24, 23).Value2 = 1.7  # W24: 1.67 -> 1.7
